# Bauble-track section relationship has been refactored and neatened.
#
# On the "Backlog" sheet, several "Bauble should notice..." / curve related
# task rows are now marked complete (column C = "X"), the now-superseded
# task in row 171 is struck through instead, and the view is scrolled/split
# to show the newly-updated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Mark rows 165-170 and 172-173 as complete in column C ("X" = done).
$ws.Range("C165").Value = "X"
$ws.Range("C166").Value = "X"
$ws.Range("C167").Value = "X"
$ws.Range("C168").Value = "X"
$ws.Range("C169").Value = "X"
$ws.Range("C170").Value = "X"
$ws.Range("C172").Value = "X"
$ws.Range("C173").Value = "X"

# Row 171's task is now obsolete/superseded -> struck-through instead of
# being checked off in column C.
$ws.Range("B171").Font.Strikethrough = $true

# Reflect the updated view/scroll state: the split stays between rows
# 164/165 but the window has scrolled down so row 165 is the first visible
# row under the split, with the last active selection on C170.
$win = $excel.ActiveWindow
$win.SplitRow = 164
$ws.Range("A164").Select() | Out-Null
$ws.Range("C170").Select() | Out-Null
